$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("B4").Value = 'Rincón De Romos'
$ws.Range("A23").Value = 'Ciudad De México'
$ws.Range("A32").Value = 'Estado De México'
$ws.Range("B32").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B33").Value = 'Almoloya De Juárez'
$ws.Range("B43").Value = 'Apaseo El Alto'
$ws.Range("B51").Value = 'Purísima Del Rincón'
$ws.Range("B53").Value = 'San Francisco Del Rincón'
$ws.Range("B54").Value = 'Santa Cruz De Juventino Rosas'
$ws.Range("B56").Value = 'Valle De Santiago'
$ws.Range("B58").Value = 'Acapulco De Juárez'
$ws.Range("B60").Value = 'Coyuca De Catalán'
$ws.Range("B62").Value = 'Zihuatanejo De Azueta'
$ws.Range("B64").Value = 'Técpan De Galeana'
$ws.Range("D66").Value = 0.0964467005076142
$ws.Range("B68").Value = 'Cuautepec De Hinojosa'
$ws.Range("B71").Value = 'Tulancingo De Bravo'
$ws.Range("B73").Value = 'Atotonilco El Alto'
$ws.Range("B79").Value = 'San Diego De Alejandría'
$ws.Range("B81").Value = 'Tepatitlán De Morelos'
$ws.Range("B84").Value = 'Unión De San Antonio'
$ws.Range("D86").Value = 0.0964467005076142
$ws.Range("B102").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B123").Value = 'Villa De Ramos'
$ws.Range("B130").Value = 'Nacozari De García'
$ws.Range("B136").Value = 'Sanctórum De Lázaro Cárdenas'
$ws.Range("B140").Value = 'Cosautlán De Carvajal'
$ws.Range("B142").Value = 'Ignacio De La Llave'
$ws.Range("B143").Value = 'Ixhuatlán Del Café'
$ws.Range("B144").Value = 'Martínez De La Torre'
$ws.Range("B145").Value = 'Medellín De Bravo'
$ws.Range("B148").Value = 'Soledad De Doblado'

$ws.Rows("159:163").Delete()

